# Update division problems in the table to the new values.
$d = $word.ActiveDocument

$replacements = @(
    @("68÷4=", "80÷8="),
    @("64÷5=", "68÷6="),
    @("49÷6=", "12÷5="),
    @("88÷5=", "64÷7="),
    @("30÷7=", "76÷2="),
    @("32÷4=", "17÷7="),
    @("69÷4=", "79÷6="),
    @("40÷6=", "39÷3="),
    @("38÷8=", "91÷7="),
    @("64÷3=", "91÷2="),
    @("98÷8=", "86÷9="),
    @("75÷4=", "25÷8="),
    @("77÷7=", "12÷8="),
    @("27÷6=", "38÷5="),
    @("93÷2=", "60÷8="),
    @("51÷5=", "99÷5="),
    @("92÷6=", "30÷4="),
    @("93÷6=", "85÷6="),
    @("43÷7=", "98÷6="),
    @("85÷8=", "15÷5="),
    @("35÷8=", "87÷8="),
    @("88÷8=", "69÷7="),
    @("89÷2=", "72÷8="),
    @("74÷4=", "56÷6="),
    @("67÷5=", "81÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
